# Closing Review points 2, 8
# Update the "Status" column (H) for CR_REVIEW_002 (row 3) and
# CR_REVIEW_008 (row 9) from "Open" to "Closed", and leave the sheet
# scrolled/zoomed/selected the way it was left after closing those points.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H3").Value = "Closed"
$ws.Range("H9").Value = "Closed"

$ws.Select()
$ws.Range("H9").Select()
$excel.ActiveWindow.Zoom = 85
